$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Establish the external workbook reference -----------------------------
# Writing a formula that points at an external workbook (bracket syntax)
# forces the engine to materialise the xl/externalLinks/externalLink1.xml
# part + the workbook-level <externalReferences> entry. We do this in a
# scratch cell, well outside the sheet's real used range, then clear it.
$scratch = $ws.Range("Z100")
$scratch.Formula = "='[Book1.xlsx]sheet'!A1"
$scratch.ClearContents()

# --- Defined names -----------------------------------------------------
# BROKEN references the external workbook/sheet registered above (now
# resolvable through the numeric [1] external-book index).
$wb.Names.Add('BROKEN', '=[1]sheet!XFC1048576')

# OUTPUT points at the new total cell on the DATA sheet.
$wb.Names.Add('OUTPUT', '=DATA!$B$5')

# --- New formula row -----------------------------------------------------
$ws.Range("B5").Formula = '=SUM(B2,B4,IFERROR(BROKEN,0))'

# --- Selection update ------------------------------------------------------
$null = $ws.Range("B6").Select()
